$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sample user row with new test values
$ws.Range("A2").Value = "testuser1"
$ws.Range("B2").Value = "test1"
$ws.Range("C2").Value = "user1"

# Extend the Email column's hyperlink-style formatting down to new filter rows
$ws.Range("G2").Copy()
$ws.Range("G3:G8").PasteSpecial(-4122)

# Register the small (8pt) helper font used for the sheet's phonetic guide info
$tempStyle = $wb.Styles.Add("PhoneticGuide")
$tempStyle.Font.Size = 8
$tempStyle.Delete()

[void]$ws.Range("J16").Select()
